$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 132
$ws1.Range("F8").Value = 4814
$ws1.Range("F13").Value = 93

# Sheet "全部类型" (All types) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 132
$ws4.Range("F9").Value = 4814
$ws4.Range("F14").Value = 93
